$d = $word.ActiveDocument

# The third table in the document (3 columns x 2 rows) currently has an
# "auto" preferred width and a third column that is 3645 twips (182.25 pt)
# wide. Widen that column to 3933 twips (196.65 pt) and switch the table's
# preferred width to a fixed value (8651 twips / 432.55 pt), matching the
# sum of the (unchanged) first two columns plus the new third column width.
$t = $d.Tables.Item(3)

$t.Columns.Item(3).Width = 196.65

$t.PreferredWidthType = 3
$t.PreferredWidth = 432.55
